# EffectData.xlsx - "add language id for property name"
#
# The header row (row 9) holds free-text descriptions for each attribute
# column (originally localized Chinese strings). This edit replaces those
# descriptions with a stable "language id" token, built as LPID_<ID>, where
# <ID> is the machine-readable column identifier already stored in row 1.
# e.g. B1="SUCKBLOOD" / old B9="吸血"  ->  new B9="LPID_SUCKBLOOD"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 30   # column AD

for ($c = 2; $c -le $lastCol; $c++) {
    $id = $ws.Cells.Item(1, $c).Value2
    $ws.Cells.Item(9, $c).Value = "LPID_" + $id
}

# --- Best-effort view-state touch-up to mirror the author's resave ---
# (column resize + reselecting AD9 as the active cell on the frozen pane)

$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 15.571428571428571
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 8.571428571428571
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 11
$ws.Columns.Item(10).ColumnWidth = 14.571428571428571
$ws.Range("L1:Z1").EntireColumn.ColumnWidth = 15.857142857142858
$ws.Columns.Item(27).ColumnWidth = 9.571428571428571
$ws.Columns.Item(28).ColumnWidth = 12.285714285714286
$ws.Range("AC1:AD1").EntireColumn.ColumnWidth = 11.428571428571429

$ws.Range("AD9").Select()
